$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""
$ws.Range("H39").Value = 164.71428
$ws.Range("I39").Value = 164.71428
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 494.14284
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -198.14284
$ws.Range("N39").Value = ""
$ws.Range("H55").Value = 486.4
$ws.Range("I55").Value = 376.66666
$ws.Range("K55").Value = 376.66666
$ws.Range("M55").Value = -162.66666
$ws.Range("H58").Value = 930.93335
$ws.Range("I58").Value = 542.1818
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1626.5454
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -1476.5454
$ws.Range("N58").Value = -6300
$ws.Range("H106").Value = 4558
$ws.Range("I106").Value = 4558
$ws.Range("K106").Value = 4558
$ws.Range("M106").Value = -3927
$ws.Range("H116").Value = 3676
$ws.Range("I116").Value = 2992.5
$ws.Range("K116").Value = 2992.5
$ws.Range("M116").Value = 449.5
$ws.Range("H125").Value = 3840.6155
$ws.Range("I125").Value = 2193
$ws.Range("J125").Value = 9332.666999999999
$ws.Range("K125").Value = 19737
$ws.Range("L125").Value = 83994.003
$ws.Range("M125").Value = -17277
$ws.Range("N125").Value = -88914.003
$ws.Range("H138").Value = 3414.7727
$ws.Range("J138").Value = 4878.727
$ws.Range("L138").Value = 14636.181
$ws.Range("N138").Value = -24916.181
$ws.Range("H141").Value = 3293.2144
$ws.Range("I141").Value = 3486.1538
$ws.Range("K141").Value = 10458.4614
$ws.Range("M141").Value = -5278.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4995.8887
$ws.Range("I32").Value = 4995.8887
$ws.Range("K32").Value = 4995.8887
$ws.Range("M32").Value = -4708.8887
$ws.Range("H102").Value = 2041.8182
$ws.Range("I102").Value = 1382
$ws.Range("J102").Value = 3801.3333
$ws.Range("K102").Value = 1382
$ws.Range("L102").Value = 3801.3333
$ws.Range("M102").Value = 240
$ws.Range("N102").Value = -7045.3333
$ws.Range("H132").Value = 1757.0769
$ws.Range("I132").Value = 1655.5555
$ws.Range("K132").Value = 4966.666499999999
$ws.Range("M132").Value = -2436.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3030.72
$ws.Range("I20").Value = 2327.1875
$ws.Range("K20").Value = 2327.1875
$ws.Range("M20").Value = -2080.1875
$ws.Range("H86").Value = 2455.8
$ws.Range("I86").Value = 2455.8
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2455.8
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1332.8
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 2455.8
$ws.Range("I89").Value = 2455.8
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 12279
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -6663
$ws.Range("N89").Value = ""
$ws.Range("H105").Value = 1713.2858
$ws.Range("I105").Value = 1268.6
$ws.Range("K105").Value = 1268.6
$ws.Range("M105").Value = 478.4000000000001
$ws.Range("H134").Value = 6003.6924
$ws.Range("I134").Value = 6604.8
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 19814.4
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -17279.4
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2346.9285
$ws.Range("I31").Value = 1832.1333
$ws.Range("J31").Value = 3633.9167
$ws.Range("K31").Value = 1832.1333
$ws.Range("L31").Value = 3633.9167
$ws.Range("M31").Value = -1537.1333
$ws.Range("N31").Value = -4223.9167
$ws.Range("H34").Value = 2346.9285
$ws.Range("I34").Value = 1832.1333
$ws.Range("J34").Value = 3633.9167
$ws.Range("K34").Value = 1832.1333
$ws.Range("L34").Value = 3633.9167
$ws.Range("M34").Value = -1630.1333
$ws.Range("N34").Value = -4037.9167
$ws.Range("H47").Value = 15000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 15000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 15000
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = -16132
$ws.Range("H50").Value = 27628.125
$ws.Range("J50").Value = 28563.143
$ws.Range("L50").Value = 28563.143
$ws.Range("N50").Value = -29813.143
$ws.Range("H54").Value = 9810.5
$ws.Range("J54").Value = 9810.5
$ws.Range("L54").Value = 9810.5
$ws.Range("N54").Value = -11126.5
$ws.Range("H60").Value = 17666.666
$ws.Range("H132").Value = 2374.75
$ws.Range("I132").Value = 1999
$ws.Range("K132").Value = 5997
$ws.Range("M132").Value = -3467
$ws.Range("H134").Value = 948.5294
$ws.Range("I134").Value = 958.3570999999999
$ws.Range("K134").Value = 2875.0713
$ws.Range("M134").Value = -340.0712999999996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 14825.25
$ws.Range("I26").Value = 4500.5
$ws.Range("J26").Value = 25150
$ws.Range("K26").Value = 13501.5
$ws.Range("L26").Value = 75450
$ws.Range("M26").Value = -13213.5
$ws.Range("N26").Value = -76026
$ws.Range("H41").Value = 2000
$ws.Range("J41").Value = 2000
$ws.Range("L41").Value = 6000
$ws.Range("N41").Value = -6676
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H114").Value = 616.3333
$ws.Range("I114").Value = 599.5
$ws.Range("J114").Value = 624.75
$ws.Range("K114").Value = 1798.5
$ws.Range("L114").Value = 1874.25
$ws.Range("M114").Value = 1455.5
$ws.Range("N114").Value = -8382.25
$ws.Range("H129").Value = 982
$ws.Range("I129").Value = 467.5
$ws.Range("K129").Value = 1402.5
$ws.Range("M129").Value = 3597.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 7459.8
$ws.Range("I126").Value = 8249.75
$ws.Range("J126").Value = 4300
$ws.Range("K126").Value = 24749.25
$ws.Range("L126").Value = 12900
$ws.Range("M126").Value = -22279.25
$ws.Range("N126").Value = -17840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 983.9375
$ws.Range("I22").Value = 1032.5555
$ws.Range("J22").Value = 921.4286
$ws.Range("K22").Value = 1032.5555
$ws.Range("L22").Value = 921.4286
$ws.Range("M22").Value = -737.5554999999999
$ws.Range("N22").Value = -1511.4286
$ws.Range("H27").Value = 983.9375
$ws.Range("I27").Value = 1032.5555
$ws.Range("J27").Value = 921.4286
$ws.Range("K27").Value = 1032.5555
$ws.Range("L27").Value = 921.4286
$ws.Range("M27").Value = -925.5554999999999
$ws.Range("N27").Value = -1135.4286
$ws.Range("H46").Value = 3823.4707
$ws.Range("I46").Value = 2833.3333
$ws.Range("J46").Value = 4363.5454
$ws.Range("K46").Value = 2833.3333
$ws.Range("L46").Value = 4363.5454
$ws.Range("M46").Value = -2645.3333
$ws.Range("N46").Value = -4739.5454
$ws.Range("H74").Value = 49997
$ws.Range("J74").Value = 49997
$ws.Range("L74").Value = 49997
$ws.Range("N74").Value = -51993
$ws.Range("H77").Value = 49997
$ws.Range("J77").Value = 49997
$ws.Range("L77").Value = 149991
$ws.Range("N77").Value = -159975
$ws.Range("H132").Value = 8158.5835
$ws.Range("J132").Value = 8749
$ws.Range("L132").Value = 26247
$ws.Range("N132").Value = -31307
$ws.Range("H136").Value = 3000.8572
$ws.Range("I136").Value = 2834.3333
$ws.Range("K136").Value = 8502.999899999999
$ws.Range("M136").Value = -5952.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 6700
$ws.Range("I58").Value = 6700
$ws.Range("K58").Value = 6700
$ws.Range("M58").Value = -6392
$ws.Range("H132").Value = 1523.7646
$ws.Range("I132").Value = 1557.25
$ws.Range("K132").Value = 4671.75
$ws.Range("M132").Value = -2141.75
